# Apply the "Add specifications for DQ rules for base_product" edit:
#  - add a note cell on Sheet1
#  - add a new "Base_Product" worksheet after Sheet1
#  - populate it with the new DQ-rule specification table (Table13)
#  - format a few columns/cells to match the authored layout
#  - make Base_Product the active/selected sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: small addendum note next to the CON.2 boundary row ---
$ws1.Activate()
$ws1.Range("I11").Value = "e.g.: minum effect_date?"
$ws1.Range("H12").Select()

# --- Create the new Base_Product worksheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Base_Product"

# --- Header row ---
$ws2.Range("B2").Value = "Rule_id"
$ws2.Range("C2").Value = "DataCategory"
$ws2.Range("D2").Value = "Element"
$ws2.Range("E2").Value = "Enforced?"
$ws2.Range("F2").Value = "Allowed_values"
$ws2.Range("G2").Value = "Lower_limit"
$ws2.Range("H2").Value = "Upper_limit"
$ws2.Range("I2").Value = "Replacement_value"
$ws2.Range("J2").Value = "Required"

# --- Data rows ---
# CSV rows
$ws2.Range("B3").Value = "CSV.1"
$ws2.Range("C3").Value = "Demographic"
$ws2.Range("E3").Value = "Y"

$ws2.Range("B4").Value = "CSV.2"
$ws2.Range("C4").Value = "Demographic"
$ws2.Range("E4").Value = "N"

$ws2.Range("B5").Value = "CSV.3"
$ws2.Range("C5").Value = "Demographic"
$ws2.Range("E5").Value = "Y"

$ws2.Range("B6").Value = "CSV.4"
$ws2.Range("C6").Value = "Demographic"
$ws2.Range("E6").Value = "N"
$ws2.Range("G6").Value = 10000

$ws2.Range("B7").Value = "CSV.5"
$ws2.Range("C7").Value = "Demographic"
$ws2.Range("E7").Value = "Y"

# SCH / Schema rows
$ws2.Range("B8").Value = "SCH.1"
$ws2.Range("C8").Value = "Demographic"
$ws2.Range("E8").Value = "Y"

$ws2.Range("B9").Value = "SCH.2"
$ws2.Range("C9").Value = "Demographic"
$ws2.Range("E9").Value = "Y"

# BDY.2 - credit score boundary
$ws2.Range("B10").Value = "BDY.2"
$ws2.Range("C10").Value = "Demographic"
$ws2.Range("D10").Value = "CREDIT_SCORE"
$ws2.Range("G10").Value = 300
$ws2.Range("H10").Value = 850
$ws2.Range("I10").Value = 550

# CON.1 - credit score source
$ws2.Range("B11").Value = "CON.1"
$ws2.Range("C11").Value = "Demographic"
$ws2.Range("D11").Value = "CREDIT_SCORE_SOURCE"
$ws2.Range("F11").Value = '"Equifax", "Experian", "TransUnion", NULL'

# BDY.2 - effective date boundary
$ws2.Range("B12").Value = "BDY.2"
$ws2.Range("C12").Value = "Demographic"
$ws2.Range("D12").Value = "EFFECTIVE_DATE"
$ws2.Range("G12").Value = "lower_limit_date"
$ws2.Range("H12").Value = "upper_limit_date"
$ws2.Range("I12").Value = "?"

# SCH.1 - Filmographic / Industry codes
$ws2.Range("B13").Value = "SCH.1"
$ws2.Range("C13").Value = "Filmographic"
$ws2.Range("D13").Value = "INDUSTRTY_NAICS"
$ws2.Range("J13").Value = "Required if Industry_SIC  not provided"

$ws2.Range("B14").Value = "SCH.1"
$ws2.Range("C14").Value = "Filmographic"
$ws2.Range("D14").Value = "INDUSTRY_SIC"
$ws2.Range("J14").Value = "Required if Industry_NAICS  not provided"

# --- Formatting ---
$ws2.Range("C3:C9").WrapText = $true
$ws2.Range("C11").WrapText = $true
$ws2.Range("C13:C14").WrapText = $true
$ws2.Range("C15:C22").WrapText = $true

$ws2.Range("I2").HorizontalAlignment = -4131   # xlLeft
$ws2.Range("I10").HorizontalAlignment = -4152  # xlRight
$ws2.Range("I12").HorizontalAlignment = -4152  # xlRight

$ws2.Range("G11:G22").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("B15:B22").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("H11").NumberFormat = "General"
$ws2.Range("H18:H22").NumberFormat = "General"

# --- Column widths (approximate the authored layout) ---
$ws2.Columns.Item(2).ColumnWidth = 8.86328125
$ws2.Columns.Item(3).ColumnWidth = 22.6640625
$ws2.Columns.Item(4).ColumnWidth = 17.33203125
$ws2.Columns.Item(5).ColumnWidth = 17.33203125
$ws2.Columns.Item(6).ColumnWidth = 34.46484375
$ws2.Columns.Item(7).ColumnWidth = 14.265625
$ws2.Columns.Item(8).ColumnWidth = 14.3984375
$ws2.Columns.Item(9).ColumnWidth = 19.06640625
$ws2.Columns.Item(10).ColumnWidth = 33.265625

# --- Turn the range into a table (ListObject) ---
$lo = $ws2.ListObjects.Add(1, $ws2.Range("B2:J22"), $null, 1)
$lo.Name = "Table13"
$lo.TableStyle = "TableStyleMedium4"

# --- Make Base_Product the active / selected sheet ---
$ws2.Activate()
$ws2.Range("G6").Select()
